# Edit script: remove the "by Hines, Raghav, and Khan" attribution from the
# introduction paragraph (commit: "Removed previous author from word documents"),
# relocate the Word "_GoBack" last-edit bookmark to land inside the word
# "document" (splitting it into "doc" / "ument"), and split the "Agile" run
# in the following paragraph the same way Word's live grammar checker does
# (it wraps the word in its own run so it can attach proofing markers).
#
# The bookmark name "_GoBack" is unique per document, so re-adding it at the
# new location automatically removes it from its old location (next to
# "... MVC Core") -- matching the third hunk of the diff.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Introduction paragraph: drop " by Hines, Raghav, and Khan" and leave
#    the "_GoBack" bookmark sitting between "doc" and "ument" of the word
#    "document" (this also forces that word to be split into two runs).
# ---------------------------------------------------------------------

# Add a throw-away bookmark right at the boundary that precedes the
# existing " classe" run. Find/Replace's text edit normalises (merges)
# same-formatted runs within the paragraph it touches; anchoring a
# bookmark at that boundary first keeps the later runs (" classe", "s and
# controllers...", "wireframe diagrams...") from being swept into the
# merge.
$guard = $d.Content
$guard.Find.Execute("model classe")
$guardPos = $guard.Start + 5   # just after "model", right before " classe"
$d.Bookmarks.Add("ZZGUARD", $d.Range($guardPos, $guardPos))

# Remove the attribution clause.
$target = $d.Content
$target.Find.Execute(" by Hines, Raghav, and Khan")
$d.Range($target.Start, $target.End).Text = ""

# Re-locate "...investigation document" now that the text has shifted, and
# remember the offsets right after "doc" and right after "document".
$locate = $d.Content
$locate.Find.Execute("The author has examined the initial investigation doc")
$afterDoc = $locate.End
$afterDocument = $afterDoc + 5

# Force the run split between "doc" and "ument" (a temporary bookmark,
# added then immediately removed, leaves the run boundary behind without
# leaving any bookmark markup there).
$d.Bookmarks.Add("ZZSPLIT", $d.Range($afterDoc, $afterDoc))
$d.Bookmarks("ZZSPLIT").Delete()

# Drop the real "_GoBack" bookmark right after "ument" -- re-adding the
# bookmark under its reserved name automatically relocates it from
# wherever it used to be (next to "... MVC Core" later in the document).
$d.Bookmarks.Add("_GoBack", $d.Range($afterDocument, $afterDocument))

# Remove the guard bookmark now that the paragraph will not be touched
# again; this leaves the existing run boundary it was protecting intact.
$d.Bookmarks("ZZGUARD").Delete()

# ---------------------------------------------------------------------
# 2) "Agile practices ..." paragraph: split "Agile" into its own run, the
#    same way Word does when the grammar checker flags it (wrapping it
#    with proofErr markers around a dedicated run).
# ---------------------------------------------------------------------

$agile = $d.Content
$agile.Find.Execute("Agile")
$agileStart = $agile.Start
$agileEnd = $agile.End

$d.Bookmarks.Add("ZZAGILE1", $d.Range($agileStart, $agileStart))
$d.Bookmarks("ZZAGILE1").Delete()

$d.Bookmarks.Add("ZZAGILE2", $d.Range($agileEnd, $agileEnd))
$d.Bookmarks("ZZAGILE2").Delete()

Write-Output "edit complete"
